$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 60-61; this shifts the existing rows 60-159 down to 62-161
# (preserving formatting/number-format from the rows being pushed down, matching
# Excel's default "insert" behaviour of carrying the format of the row above).
$ws.Rows("60:61").Insert()

# Row 60: new "Primera" quality entry for the latest week
$ws.Range("A60").Value = 11
$ws.Range("B60").Value = "Vega Monumental Concepción"
$ws.Range("C60").Value = "Bíobío"
$ws.Range("D60").Value = 44771
$ws.Range("E60").Value = 8
$ws.Range("F60").Value = 100112044
$ws.Range("G60").Value = "Perejil"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 200
$ws.Range("K60").Value = 700
$ws.Range("L60").Value = 800
$ws.Range("M60").Value = 750
$ws.Range("N60").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O60").Value = "Región de Ñuble"
$ws.Range("P60").Value = 750
$ws.Range("Q60").Value = 1
$ws.Range("R60").Value = "Hortaliza"

# Row 61: new "Segunda" quality entry for the latest week
$ws.Range("A61").Value = 11
$ws.Range("B61").Value = "Vega Monumental Concepción"
$ws.Range("C61").Value = "Bíobío"
$ws.Range("D61").Value = 44771
$ws.Range("E61").Value = 8
$ws.Range("F61").Value = 100112044
$ws.Range("G61").Value = "Perejil"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Segunda"
$ws.Range("J61").Value = 100
$ws.Range("K61").Value = 600
$ws.Range("L61").Value = 600
$ws.Range("M61").Value = 600
$ws.Range("N61").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O61").Value = "Región de Ñuble"
$ws.Range("P61").Value = 600
$ws.Range("Q61").Value = 1
$ws.Range("R61").Value = "Hortaliza"
